# Manual testcases of sprint114
# Insert a new testcase row (new row 11) above the existing "Click the
# Outlet-> Create new list" row, pushing all the later rows down by one.
# The new row documents the "Click the Checkbox" / "Clear and Download
# inventory report" test case.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 11 - everything from old row 11 downward
# (old rows 11-18 and the trailing placeholder row 22) shifts down by one
# (new rows 12-19 and 23).
$ws.Rows("11:11").Insert()

# Pick up the formatting (styles, fonts, borders) of the row immediately
# below (the row that used to be row 10's neighbour, now row 12 - i.e. the
# original row 11 content/format) so the new row matches its siblings.
$ws.Range("A12:G12").Copy()
$ws.Range("A11:G11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row height for the new row.
$ws.Rows("11:11").RowHeight = 60

# Fill in the new row's content.
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Buyer Home page"
$ws.Range("C11").Value = "View Inventory menu page"
$ws.Range("D11").Value = "Click the Checkbox "

$ws.Range("E11").Value = "Once click the Checkbox , a button will appear at the top to allow users to  ""Clear and Download inventory report"""
$chars = $ws.Range("E11").Characters(79, 38)
$chars.Font.Bold = $true

# F11/G11 stay blank (just formatted) for this new row.
$ws.Range("F11").ClearContents()
$ws.Range("G11").ClearContents()

# The SL. No column (A) is a simple running count, 1..16; renumber the rows
# that got pushed down one position so the count stays sequential.
$ws.Range("A12").Value = 11
$ws.Range("A13").Value = 12
$ws.Range("A14").Value = 13
$ws.Range("A15").Value = 14
$ws.Range("A16").Value = 15
$ws.Range("A17").Value = 16

# Restore the selection/scroll position to roughly where the user ended up
# after adding the new row (around the newly shifted rows).
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A18").Select()
